# "Cleaned up the project" - adds new AngularJS knowledge rows to the
# AngularJS sheet and tweaks a few existing row heights.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AngularJS")
$ws.Activate()

# --- Row height tweaks on existing rows -------------------------------

# Row 17 loses its explicit 30pt height (back to default row height).
$ws.Rows.Item(17).AutoFit() | Out-Null

# Row 34: 60 -> 45
$ws.Rows.Item(34).RowHeight = 45

# Row 40: 90 -> 75
$ws.Rows.Item(40).RowHeight = 75

# --- New rows 45-49 -----------------------------------------------------

$A45 = 'AngularJS Form Custom Validation'
$B45 = @'
app.directive('myDirective', function() {
  return {
    require: 'ngModel',
    link: function(scope, element, attr, mCtrl) {
      function myValidation(value) {
        if (value.indexOf("e") > -1) {
          //Add custom logic to set the value of mCtrl.$setValidity(); 
        return value;
      }
      mCtrl.$parsers.push(myValidation);
    }
  };
});
'@

$A46 = 'AngularJS Global API'
$B46 = @'
Set of global JavaScript functions for performing common tasks like:
Comparing objects
Iterating objects
Converting data
'@

$A47 = 'AngularJS Include'
$B47 = @'
Include external html file
<div ng-include="'myFile.htm'"></div>
'@

$A48 = 'Include Cross Domains'
$B48 = @'
By default, the ng-include directive does not allow to include files from other domains.
To include files from another domain, we can add a whitelist of legal files and/or domains in the config function of the application
app.config(function($sceDelegateProvider) {
    $sceDelegateProvider.resourceUrlWhitelist([
        'https://tryit.w3schools.com/**'
    ]);
});
'@

$A49 = 'AngularJS Routing'
$B49 = 'If we want to navigate to different pages in the application, but also want the application to be a SPA (Single Page Application), with no page reloading, we can use the ngRoute module.'

$ws.Range("A45").Value = $A45
$ws.Range("A45").Style = "Normal"
$ws.Range("B45").Value = $B45
$ws.Range("B45").WrapText = $true

$ws.Range("A46").Value = $A46
$ws.Range("A46").Style = "Normal"
$ws.Range("B46").Value = $B46
$ws.Range("B46").WrapText = $true

$ws.Range("A47").Value = $A47
$ws.Range("A47").Style = "Normal"
$ws.Range("B47").Value = $B47
$ws.Range("B47").WrapText = $true

# NOTE: shared-string insertion order in the source workbook has both
# row-48 and row-49 labels (column A) entered before either body text
# (column B), so we replicate that exact order here.
$ws.Range("A48").Value = $A48
$ws.Range("A48").Style = "Normal"
$ws.Range("A49").Value = $A49
$ws.Range("A49").Style = "Normal"
$ws.Range("B48").Value = $B48
$ws.Range("B48").WrapText = $true
$ws.Range("B49").Value = $B49
$ws.Range("B49").WrapText = $true

$ws.Rows.Item(45).RowHeight = 195
$ws.Rows.Item(46).RowHeight = 60
$ws.Rows.Item(47).RowHeight = 30
$ws.Rows.Item(48).RowHeight = 120
$ws.Rows.Item(49).RowHeight = 30

# --- View state: scroll frozen pane + selection --------------------------

$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B54").Select()
